$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting for the new year columns (E:H) from column D ---
$ws.Range("D3").Copy()
$ws.Range("E3:H3").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("E4:H4").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("E5:H5").PasteSpecial(-4122)

$ws.Range("D6").Copy()
$ws.Range("E6:H6").PasteSpecial(-4122)

# --- Row 3: years ---
$ws.Range("E3").Value() = 2020
$ws.Range("F3").Value() = 2021
$ws.Range("G3").Value() = 2022
$ws.Range("H3").Value() = 2023

# --- Row 4: number of local governments (484 every year, now numeric) ---
$ws.Range("D4").Value() = 484
$ws.Range("E4").Value() = 484
$ws.Range("F4").Value() = 484
$ws.Range("G4").Value() = 484
$ws.Range("H4").Value() = 484

# --- Row 5: proportion values ---
$ws.Range("E5").Value() = 13.2
$ws.Range("F5").Value() = 21.5
$ws.Range("G5").Value() = 34.5
$ws.Range("H5").Value() = 40.53

# --- Row 6: counts ---
$ws.Range("E6").Value() = 67
$ws.Range("F6").Value() = 104
$ws.Range("G6").Value() = 167
$ws.Range("H6").Value() = 169

# --- Header row resizing (text/content is unchanged, only row height & col width) ---
$ws.Rows(1).RowHeight = 79.5
$ws.Range("A1:C1").ColumnWidth = 39.33

# --- Remove the stray selection left over from the editing session ---
$ws.Range("A1").Select()
